$d = $word.ActiveDocument

# Locate the bullet paragraph "Do some early life conditions affect an
# adult individual's metabolic status and risk for disease more than other
# conditions?" - the new bullet about long-term diabetes outcomes goes
# right after it (and right before the "To what extent does the family
# network ..." bullet). We resolve it to its 1-based position in
# $d.Paragraphs up front (Paragraph.Index is not a reliable 1:1 mapping to
# that collection's ordinal position in this host, so we track the
# position ourselves instead of recomputing it from .Index later).
$needle = "status and risk for disease more than other conditions?"
$anchorPos = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Contains($needle)) {
        $anchorPos = $i
        break
    }
}
if ($anchorPos -eq 0) {
    throw "Anchor paragraph not found"
}

# Collapse to the end of that paragraph's text and add a new paragraph
# after it. InsertParagraphAfter carries over the pPr (pStyle "Compact" +
# the numPr ilvl=0/numId=1001 bullet numbering) from the preceding
# paragraph, so the new bullet automatically becomes part of the same
# list.
$anchorRange = $d.Paragraphs($anchorPos).Range
$anchorRange.Collapse(0) | Out-Null
$anchorRange.InsertParagraphAfter() | Out-Null

$newParaPos = $anchorPos + 1
$insertAt = $d.Paragraphs($newParaPos).Range.Start

# Build the new bullet's text as five separate runs (three text segments
# split by single-space runs), matching how the rest of the document's
# wrapped paragraphs are authored. Each segment is typed into its own
# mini-paragraph and the paragraph marks between them are then deleted -
# this keeps every segment in its own <w:r> run instead of the run-merging
# that happens when text is typed straight into one run via repeated
# InsertAfter calls.
$segments = @(
    "Are long-term diabetes outcomes (nephropathy, eye disease, premature",
    " ",
    "death) affected by treatment persistence and stability in patients",
    " ",
    "with psychiatric comorbidity?"
)

$work = $d.Range($insertAt, $insertAt)
for ($i = 0; $i -lt $segments.Length; $i++) {
    $work.InsertAfter($segments[$i]) | Out-Null
    if ($i -lt $segments.Length - 1) {
        $work.Collapse(0) | Out-Null
        $work.InsertParagraphAfter() | Out-Null
        $work.Collapse(0) | Out-Null
        $work.MoveStart(1, 1) | Out-Null
    }
}

# Join the mini-paragraphs back into a single paragraph by deleting the
# paragraph marks that separate them, which preserves each segment as its
# own run instead of coalescing them.
for ($i = 0; $i -lt $segments.Length - 1; $i++) {
    $p = $d.Paragraphs($newParaPos)
    $mark = $d.Range($p.Range.End - 1, $p.Range.End)
    $mark.Delete() | Out-Null
}

Write-Output "New paragraph text: $($d.Paragraphs($newParaPos).Range.Text)"
